$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.130.71'
$ws.Range("E2").Value = '  -0.95%  '
$ws.Range("D3").Value = '2.945.32'
$ws.Range("E3").Value = '  -1.74%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '''375.17'
$ws.Range("E5").Value = '  -1.69%  '
$ws.Range("D6").Value = '''101.38'
$ws.Range("E6").Value = '  -3.59%  '
$ws.Range("E7").Value = '  -1.86%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").Value = '''0.586'
$ws.Range("E9").Value = '  -2.25%  '
$ws.Range("D10").Value = '''36.40'
$ws.Range("E10").Value = '  -3.10%  '
$ws.Range("E11").Value = '  -0.91%  '
$ws.Range("D12").Value = '''0.0851'
$ws.Range("E12").Value = '  +0.01%  '
$ws.Range("D13").Value = '3.412.02'
$ws.Range("E13").Value = '  -1.35%  '
$ws.Range("D14").Value = '''18.07'
$ws.Range("D15").Value = '''7.57'
$ws.Range("E15").Value = '  -0.71%  '
$ws.Range("D16").Value = '2.937.50'
$ws.Range("E16").Value = '  -1.95%  '
$ws.Range("E17").Value = '  +1.27%  '
$ws.Range("D18").Value = '''10.67'
$ws.Range("E18").Value = '  +42.83%  '
$ws.Range("D19").Value = '51.073.85'
$ws.Range("E19").Value = '  -1.01%  '
$ws.Range("D20").Value = '''3.10'
$ws.Range("E20").Value = '  -7.13%  '
$ws.Range("D21").Value = '''12.46'
$ws.Range("E21").Value = '  -4.58%  '
$ws.Range("E22").Value = '  -1.07%  '
$ws.Range("D23").Value = '''265.97'
$ws.Range("E23").Value = '  +0.93%  '
$ws.Range("D24").Value = '''68.72'
$ws.Range("E24").Value = '  -1.13%  '
$ws.Range("D25").Value = '''3.15'
$ws.Range("E25").Value = '  +8.31%  '
$ws.Range("D26").Value = '''8.18'
$ws.Range("E26").Value = '  -2.54%  '
$ws.Range("D27").Value = '''7.62'
$ws.Range("E27").Value = '  -2.09%  '
$ws.Range("D28").Value = '''0.999'
$ws.Range("E28").Value = '  -0.10%  '
$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").Value = '''0.164'
$ws.Range("E29").Value = '  -4.57%  '
$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").Value = '''25.65'
$ws.Range("E30").Value = '  -1.81%  '
$ws.Range("E31").Value = '  -6.05%  '
$ws.Range("D32").Value = '''10.02'
$ws.Range("E32").Value = '  +0.85%  '
$ws.Range("D33").Value = '''50.70'
$ws.Range("E33").Value = '  -0.79%  '
$ws.Range("E34").Value = '  -1.43%  '
$ws.Range("D35").Value = '''33.40'
$ws.Range("E35").Value = '  -5.37%  '
$ws.Range("D36").Value = '''0.0443'
$ws.Range("E36").Value = '  -2.84%  '
$ws.Range("E37").Value = '  -0.10%  '
$ws.Range("E38").Value = '  +2.92%  '
$ws.Range("E39").Value = '  -1.56%  '
$ws.Range("D40").Value = '''16.38'
$ws.Range("E40").Value = '  -5.33%  '
$ws.Range("E41").Value = '  -3.66%  '
$ws.Range("E42").Value = '  -4.37%  '
$ws.Range("D43").Value = '''120.58'
$ws.Range("E43").Value = '  -4.22%  '
$ws.Range("D44").Value = '''21.31'
$ws.Range("E44").Value = '  -2.64%  '
$ws.Range("E45").Value = '  -1.10%  '
$ws.Range("E46").Value = '  +1.41%  '
$ws.Range("D47").Value = '''0.273'
$ws.Range("E47").Value = '  -4.17%  '
$ws.Range("D48").Value = '''2.30'
$ws.Range("E48").Value = '  -2.90%  '
$ws.Range("D49").Value = '1.991.01'
$ws.Range("E49").Value = '  -2.70%  '
$ws.Range("D50").Value = '''0.0325'
$ws.Range("E50").Value = '  -2.51%  '
$ws.Range("D51").Value = '''1.31'
$ws.Range("E51").Value = '  +1.19%  '
